# Update the ARIMA summary-stats worksheet:
#  - a new outlier row "AO2001Sep" is inserted right after the header row
#  - a new outlier row "AO2021Jul" is inserted before the "ar.L1" block
#  - all of the model-coefficient rows are refreshed with the re-fit values
#    (X13 was adjusted to only look for LS/AO, no other regressors)
#  - the sheet grows from A1:G11 to A1:G13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final row contents, in top-to-bottom order, starting at row 2.
# Each entry: Label, coef, std err, z, P>|z|, [0.025, 0.975]
$rows = @(
    @("AO2001Sep", -443900,              98100,              -4.526,   0,      -636000,  -252000),
    @("AO2008Sep", -371000,              82100,              -4.517,   0,      -532000,  -210000),
    @("AO2020Mar", 90680,                24400,              3.716,    0,      42800,    139000),
    @("AO2021Jul", 197600,               644.773,            306.456,  0,      196000,   199000),
    @("ar.L1",     -1.0988,              0.102,              -10.755,  0,      -1.299,   -0.899),
    @("ar.L2",     -0.7786999999999999,  0.108,              -7.19,    0,      -0.991,   -0.5659999999999999),
    @("ma.L1",     1.1724,               0.08500000000000001, 13.826,  0,      1.006,    1.339),
    @("ma.L2",     0.9308,               0.08599999999999999, 10.77,   0,      0.761,    1.1),
    @("ar.S.L12",  0.9405,               0.04,               23.283,   0,      0.861,    1.02),
    @("ma.S.L12",  -0.9023,              0.073,              -12.43,   0,      -1.045,   -0.76),
    @("ma.S.L24",  0.1085,               0.052,              2.096,    0.036,  0.007,    0.21),
    @("sigma2",    14210000000,          1.814,              7830000000, 0,    14200000000, 14200000000)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}

# Rows 12 and 13 are brand new - give column A the same label style
# (bold, centered, thin border) used by the rest of the column.
$ws.Range("A2").Copy()
$ws.Range("A12:A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the labels (PasteSpecial only copies formatting, not values).
$ws.Cells.Item(12, 1).Value = "ma.S.L24"
$ws.Cells.Item(13, 1).Value = "sigma2"
